# Update crypto price/volume data cells per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.275.59"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "3.014.33"
$ws.Range("E3").Value = "  -3.51%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "567.14"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -3.07%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "128.88"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -4.18%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.013.32"
$ws.Range("E8").Value = "  -3.48%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.498"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.40%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.135"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.91%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.22"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  -4.48%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000223"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.73%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "32.83"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.10%  "
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "3.509.52"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("D17").Value = "61.296.20"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "3.013.21"
$ws.Range("E18").Value = "  -3.86%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.23"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -4.12%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "438.76"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.43%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.16"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -4.40%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.662"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -5.05%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.15"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -4.84%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "79.05"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -4.74%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "12.54"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -5.09%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("E28").Value = "  -6.18%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.19"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -5.46%  "
$ws.Range("E30").Value = "  -7.34%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "25.60"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -5.15%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.88"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -5.77%  "
$ws.Range("E33").Value = "  -8.12%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("E35").Value = "  -6.14%  "
$ws.Range("E36").Value = "  -3.43%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "50.03"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "0.0₃0677"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("E39").Value = "  -5.00%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "7.74"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.108"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "372.55"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -5.20%  "
$ws.Range("D43").Value = "2.654.52"
$ws.Range("E43").Value = "  -3.33%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.43"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -7.77%  "
$ws.Range("E46").Value = "  -4.48%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "119.47"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.81%  "
$ws.Range("E48").Value = "  -5.91%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "32.95"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -3.12%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.106"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.16%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "23.58"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -5.76%  "
